# Refresh the cryptos list: update Price (col D) and Volume(1h) (col E)
# for each coin row, and swap the Celestia / FTXToken rows (45 and 48)
# whose relative ranking changed.
#
# Leading apostrophes on column D assignments force plain numeric-looking
# strings (e.g. "231.43") to be stored as text, matching the sheet's
# original inlineStr/text cell type instead of being coerced to a Number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "43.696.61"
$ws.Cells.Item(2, 5).Value = "  +5.12%  "
# Row 3
$ws.Cells.Item(3, 4).Value = "2.271.58"
$ws.Cells.Item(3, 5).Value = "  +2.79%  "
# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.08%  "
# Row 5
$ws.Cells.Item(5, 4).Value = "'231.43"
$ws.Cells.Item(5, 5).Value = "  +1.07%  "
# Row 6
$ws.Cells.Item(6, 4).Value = "'0.628"
$ws.Cells.Item(6, 5).Value = "  +2.06%  "
# Row 7
$ws.Cells.Item(7, 4).Value = "'63.34"
$ws.Cells.Item(7, 5).Value = "  +6.46%  "
# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.05%  "
# Row 9
$ws.Cells.Item(9, 4).Value = "'0.432"
$ws.Cells.Item(9, 5).Value = "  +7.96%  "
# Row 10
$ws.Cells.Item(10, 5).Value = "  +19.17%  "
# Row 11
$ws.Cells.Item(11, 4).Value = "'57.31"
$ws.Cells.Item(11, 5).Value = "  -0.42%  "
# Row 12
$ws.Cells.Item(12, 4).Value = "'25.93"
$ws.Cells.Item(12, 5).Value = "  +16.62%  "
# Row 13
$ws.Cells.Item(13, 4).Value = "'0.104"
$ws.Cells.Item(13, 5).Value = "  +0.18%  "
# Row 14
$ws.Cells.Item(14, 4).Value = "2.607.03"
$ws.Cells.Item(14, 5).Value = "  +2.56%  "
# Row 15
$ws.Cells.Item(15, 4).Value = "'15.71"
$ws.Cells.Item(15, 5).Value = "  +2.36%  "
# Row 16
$ws.Cells.Item(16, 4).Value = "'5.94"
$ws.Cells.Item(16, 5).Value = "  +5.62%  "
# Row 17
$ws.Cells.Item(17, 4).Value = "'0.825"
$ws.Cells.Item(17, 5).Value = "  +4.40%  "
# Row 18
$ws.Cells.Item(18, 4).Value = "2.273.96"
$ws.Cells.Item(18, 5).Value = "  +2.48%  "
# Row 19
$ws.Cells.Item(19, 4).Value = "43.565.67"
$ws.Cells.Item(19, 5).Value = "  +5.01%  "
# Row 20
$ws.Cells.Item(20, 5).Value = "  +12.83%  "
# Row 21
$ws.Cells.Item(21, 4).Value = "'73.58"
$ws.Cells.Item(21, 5).Value = "  +2.65%  "
# Row 22
$ws.Cells.Item(22, 4).Value = "'6.13"
$ws.Cells.Item(22, 5).Value = "  +1.37%  "
# Row 23
$ws.Cells.Item(23, 4).Value = "'250.36"
$ws.Cells.Item(23, 5).Value = "  +3.52%  "
# Row 24
$ws.Cells.Item(24, 5).Value = "  +0.24%  "
# Row 25
$ws.Cells.Item(25, 4).Value = "'2.50"
$ws.Cells.Item(25, 5).Value = "  +6.75%  "
# Row 26
$ws.Cells.Item(26, 5).Value = "  +1.25%  "
# Row 27
$ws.Cells.Item(27, 4).Value = "'9.88"
$ws.Cells.Item(27, 5).Value = "  +2.40%  "
# Row 28
$ws.Cells.Item(28, 4).Value = "'172.26"
$ws.Cells.Item(28, 5).Value = "  +2.26%  "
# Row 29
$ws.Cells.Item(29, 4).Value = "'21.02"
$ws.Cells.Item(29, 5).Value = "  +6.74%  "
# Row 30
$ws.Cells.Item(30, 5).Value = "  -1.51%  "
# Row 31
$ws.Cells.Item(31, 4).Value = "'1.44"
$ws.Cells.Item(31, 5).Value = "  +3.03%  "
# Row 32
$ws.Cells.Item(32, 5).Value = "  +10.85%  "
# Row 33
$ws.Cells.Item(33, 5).Value = "  +1.77%  "
# Row 34
$ws.Cells.Item(34, 4).Value = "'0.0687"
$ws.Cells.Item(34, 5).Value = "  +6.32%  "
# Row 35
$ws.Cells.Item(35, 4).Value = "'5.07"
$ws.Cells.Item(35, 5).Value = "  +2.85%  "
# Row 36
$ws.Cells.Item(36, 4).Value = "'4.75"
$ws.Cells.Item(36, 5).Value = "  +2.72%  "
# Row 37
$ws.Cells.Item(37, 4).Value = "'6.82"
$ws.Cells.Item(37, 5).Value = "  +6.00%  "
# Row 38
$ws.Cells.Item(38, 4).Value = "'3.82"
$ws.Cells.Item(38, 5).Value = "  +7.59%  "
# Row 39
$ws.Cells.Item(39, 4).Value = "'2.34"
$ws.Cells.Item(39, 5).Value = "  -0.76%  "
# Row 40
$ws.Cells.Item(40, 4).Value = "'0.0249"
$ws.Cells.Item(40, 5).Value = "  +5.47%  "
# Row 41
$ws.Cells.Item(41, 5).Value = "  -0.21%  "
# Row 42
$ws.Cells.Item(42, 4).Value = "'8.40"
$ws.Cells.Item(42, 5).Value = "  -1.72%  "
# Row 43
$ws.Cells.Item(43, 4).Value = "'17.36"
$ws.Cells.Item(43, 5).Value = "  +6.12%  "
# Row 44
$ws.Cells.Item(44, 4).Value = "'0.0961"
$ws.Cells.Item(44, 5).Value = "  +0.36%  "
# Row 45
$ws.Cells.Item(45, 2).Value = "Celestia"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(45, 4).Value = "'10.36"
$ws.Cells.Item(45, 5).Value = "  +21.08%  "
# Row 46
$ws.Cells.Item(46, 5).Value = "  +0.59%  "
# Row 47
$ws.Cells.Item(47, 4).Value = "'97.70"
$ws.Cells.Item(47, 5).Value = "  +0.83%  "
# Row 48
$ws.Cells.Item(48, 2).Value = "FTXToken"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(48, 4).Value = "'4.41"
$ws.Cells.Item(48, 5).Value = "  +1.50%  "
# Row 49
$ws.Cells.Item(49, 4).Value = "1.477.24"
$ws.Cells.Item(49, 5).Value = "  +0.97%  "
# Row 50
$ws.Cells.Item(50, 4).Value = "'2.36"
$ws.Cells.Item(50, 5).Value = "  +5.69%  "
# Row 51
$ws.Cells.Item(51, 4).Value = "'1.08"
$ws.Cells.Item(51, 5).Value = "  +1.20%  "
